# which plots to use question
# Applies the semantic edit described in the commit: adds three new
# "which plot to use" related question rows and fixes up the
# "single choice" -> "schoice" labelling that already existed elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (E4): "single choice" -> "schoice"
$ws.Range("E4").Value = "schoice"

# Row 11 (E11): previously empty -> "schoice"
$ws.Range("E11").Value = "schoice"

# Row 12: new question about which plot to use
$ws.Range("A12").Value = "Understanding which plot to use "
$ws.Range("D12").Value = "penguins-which-plot-to-use"
$ws.Range("E12").Value = "schoice"

# Row 13: new question about the definition of statistic
$ws.Range("A13").Value = "Definition of statistic"
$ws.Range("D13").Value = "what-does-statistic-mean"
$ws.Range("E13").Value = "schoice"

# Row 14: new question about calculating statistics in R
$ws.Range("A14").Value = "Calculate statistics in R"
$ws.Range("D14").Value = "calculate-stuff-in-R-diamonds"
$ws.Range("E14").Value = "schoice"

# Row 16 (E16): "single choice" -> "schoice"
$ws.Range("E16").Value = "schoice"

# Move the active selection to D18 to match the saved view state
$ws.Range("D18").Select()
